# implementation of KK14 (not done)
#
# Adds a new block of observed-variable rows (77-81) to Sheet1, describing
# four more KK14 series: real private consumption/investment per capita,
# total hours worked per capita, real government expenditure per capita,
# and the labour tax rate.
#
# Columns: A = id, B = title, C = title_short, D = construction

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C77").Value = "log real private consumption per capita"
$ws.Range("D77").Value = "diff(log((PCND+ PCESV)/(index(GDPCTPI*CNP16OV))))-mean(log((PCND+ PCESV)/(index(GDPCTPI*CNP16OV))))"
$ws.Range("A78").Value = "hrw_pc_obs"
$ws.Range("B78").Value = "log total hours worked per capita"
$ws.Range("C78").Value = "log total hours worked per capita"
$ws.Range("D78").Value = "log((PRS85006023* CE16OV)/(index(CNP16OV)))-mean(log((PRS85006023* CE16OV)/(index(CNP16OV))))     "
$ws.Range("C79").Value = "log real private investment per capita"
$ws.Range("B77").Value = "change of log real private consumption per capita"
$ws.Range("B79").Value = "change of log real private investment per capita"
$ws.Range("A77").Value = "c_rgpc_obs"
$ws.Range("A79").Value = "I_rgpc_obs"
$ws.Range("D79").Value = "diff(log((GPDI+ PCDG)/(GDPCTPI*CNP16OV)))-mean(diff(log((GPDI+ PCDG)/(GDPCTPI*CNP16OV))))"
$ws.Range("B80").Value = "change of log real government expenditure per capita"
$ws.Range("C80").Value = "change of log real government expenditure per capita"
$ws.Range("A80").Value = "gexp_rgpc_obs"
$ws.Range("D80").Value = "X= (A957RC1Q027SBEA+A787RC1Q027SBEA+AD08RC1Q027SBEA-A918RC1Q027SBEA)/(GDPCTPI*CNP16OV)  ,   diff(X)-mean(diff(X))"
$ws.Range("A81").Value = "tau_w_obs"
$ws.Range("B81").Value = "log labour tax rate"
$ws.Range("C81").Value = "log labour tax rate"
$ws.Range("D81").Value = "RENTIN-CPROFIT-W255RC1Q027SBEA-PROPINC-A074RC1Q027SBEA-W071RC1Q027SBEA-WASCUR-PROPINC-COE-W780RC1Q027SBEA"

# Row 77 is a new section header/divider (like row 67): default/plain
# style, no explicit font formatting. Rows 78-81 keep the sheet's normal
# body style (inherited from the column default, same as rows 69-75).
$ws.Range("A77:D77").Style = "Normal"

# Match the saved view state: selection on C64, zoomed to 70%.
$ws.Range("C64").Select()
$excel.ActiveWindow.Zoom = 70
